$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.239.74"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "3.603.68"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.15"
$ws.Range("E5").Value = "  -2.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.25"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("D7").Value = "3.598.89"
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  +4.36%  "
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.85"
$ws.Range("E12").Value = "  -4.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000313"
$ws.Range("E13").Value = "  +7.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.65"
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").Value = "4.183.21"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.79"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "3.612.81"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "70.274.91"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.14"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.28"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.94"
$ws.Range("E24").Value = "  -7.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.26"
$ws.Range("E25").Value = "  +4.94%  "
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("E27").Value = "  -5.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.10"
$ws.Range("E28").Value = "  -4.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.39"
$ws.Range("E29").Value = "  -2.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.12"
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.61"
$ws.Range("E31").Value = "  -4.30%  "
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.85"
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("E34").Value = "  -3.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "572.48"
$ws.Range("E35").Value = "  -9.26%  "
$ws.Range("E36").Value = "  -6.25%  "
$ws.Range("D37").Value = "0.0₃0811"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  +15.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.04"
$ws.Range("E40").Value = "  +3.67%  "
$ws.Range("E41").Value = "  -4.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.54"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("E43").Value = "  -6.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.01"
$ws.Range("E44").Value = "  -5.70%  "
$ws.Range("D45").Value = "3.220.36"
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.81"
$ws.Range("E46").Value = "  +6.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0441"
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.40"
$ws.Range("E48").Value = "  +2.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.138"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  -3.81%  "
